# Ajout du lien vers le moteur de recherche des inventaires (colonne J)
# sur la feuille "Fonction" (premiere feuille du classeur).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nouvel en-tete de colonne J, avec le meme style que les autres en-tetes
# (copie de la mise en forme de I1 vers J1 puis ecriture du texte).
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "url_recherche"

# Lien (identique pour toutes les lignes de donnees) vers le moteur de
# recherche des inventaires AD13.
$searchUrl = "https://www.archives13.fr/archive/recherche/fonds/n:93"

for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 10).Value = $searchUrl
}

$excel.CutCopyMode = $false
